$d = $word.ActiveDocument

# Locate the paragraph ending in "LOM3057: Introdução aos Materiais
# Poliméricos (Requisito)" - the deletion starts right after it (removing
# the blank paragraph that follows it).
$anchor = $d.Content
$anchor.Find.Execute("LOM3057: Introdução aos Materiais Poliméricos (Requisito)") | Out-Null
$anchor.Expand(4) | Out-Null   # wdParagraph -> include the paragraph mark
$deleteStart = $anchor.End

# Locate the footer paragraph containing the copyright notice - the
# deletion ends at the end of that paragraph (inclusive of its mark).
$footer = $d.Content
$footer.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution") | Out-Null
$footer.Expand(4) | Out-Null   # wdParagraph -> include the paragraph mark
$deleteEnd = $footer.End

# Remove the blank paragraph, the "Ver no Jupiter..." paragraph and the
# copyright paragraph in one shot, leaving the subsequent blank paragraph
# and the page-break paragraph untouched.
$d.Range($deleteStart, $deleteEnd).Delete()
